$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Helper: force a full-text-range reassignment to go through even when the
# resulting text is unchanged, so the writer rebuilds the paragraph as a
# single consolidated run instead of leaving the original multi-run split.
function Set-ConsolidatedText($textRange, $finalText) {
    $textRange.Text = "~~~iron_native_tmp~~~"
    $textRange.Text = $finalText
}

# Title 1: "A" + " " + "slide" -> "A slide" (consolidate into one run)
Set-ConsolidatedText $s.Shapes.Item(1).TextFrame.TextRange "A slide"

# Table cell (row 1, col 2): "a" + " " + "table" -> "a table"
$tbl = $s.Shapes.Item(3).Table
Set-ConsolidatedText $tbl.Cell(1, 2).Shape.TextFrame.TextRange "a table"

# TextBox 3: "Plus" + " " + "an" + " " + "image" -> "Plus an image"
Set-ConsolidatedText $s.Shapes.Item(7).TextFrame.TextRange "Plus an image"
